$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "54.365.07"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "2.269.68"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "496.13"
$ws.Range("E5").Value = "  +2.26%  "
Set-TextValue "D6" "128.27"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("E7").Value = "  +0.12%  "
Set-TextValue "D8" "0.527"
$ws.Range("E8").Value = "  +1.39%  "
Set-TextValue "D9" "0.0962"
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("E10").Value = "  +2.28%  "
Set-TextValue "D11" "0.329"
$ws.Range("E11").Value = "  +4.08%  "
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "2.673.95"
$ws.Range("E13").Value = "  +2.84%  "
Set-TextValue "D14" "22.06"
$ws.Range("E14").Value = "  +4.22%  "
$ws.Range("D15").Value = "54.267.77"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "2.272.52"
$ws.Range("E17").Value = "  +1.76%  "
Set-TextValue "D18" "10.07"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "303.68"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "6.49"
$ws.Range("E21").Value = "  +5.59%  "
Set-TextValue "D22" "0.999"
$ws.Range("E22").Value = "  +0.25%  "
Set-TextValue "D23" "61.84"
$ws.Range("E23").Value = "  -1.96%  "
Set-TextValue "D24" "1.00"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "2.370.13"
$ws.Range("E25").Value = "  +2.77%  "
Set-TextValue "D26" "0.371"
$ws.Range("E26").Value = "  +1.61%  "
Set-TextValue "D27" "0.149"
$ws.Range("E27").Value = "  +1.67%  "
Set-TextValue "D28" "7.18"
$ws.Range("E28").Value = "  +1.99%  "
Set-TextValue "D29" "168.64"
$ws.Range("E29").Value = "  +3.35%  "
Set-TextValue "D30" "1.61"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").Value = "0.0₃0681"
$ws.Range("E31").Value = "  +1.75%  "
Set-TextValue "D32" "5.86"
$ws.Range("E32").Value = "  +1.80%  "
Set-TextValue "D33" "1.09"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("E34").Value = "  +0.09%  "
Set-TextValue "D35" "17.76"
$ws.Range("E35").Value = "  +2.51%  "
Set-TextValue "D36" "0.995"
$ws.Range("E36").Value = "  +0.19%  "
Set-TextValue "D37" "0.894"
$ws.Range("E37").Value = "  +6.12%  "
Set-TextValue "D38" "1.19"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  +3.64%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D40" "0.372"
$ws.Range("E40").Value = "  +1.56%  "
Set-TextValue "D41" "1.41"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "3.39"
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "126.72"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "4.77"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D45" "0.0896"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0486"
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D47" "0.546"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D48" "238.70"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0205"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D50" "10.77"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "16.20"
$ws.Range("E51").Value = "  +0.29%  "
